# Update the "enddate" column (D) for every ticker row.
# Most rows move from 2025-03-11 to 2025-03-13; a handful of rows move to
# 2025-03-14 instead.
#
# The values must remain plain text (the source workbook stores dates as
# literal strings, not Excel date serials), so each cell is briefly
# switched to Text format while the value is written, then ClearFormats()
# is used to drop that temporary formatting again so the cell's style
# stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$exceptionRows = @(162, 254, 322, 468)

$lastRow = 505

for ($r = 2; $r -le $lastRow; $r++) {
    if ($exceptionRows -contains $r) {
        $newValue = "2025-03-14"
    } else {
        $newValue = "2025-03-13"
    }

    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.ClearFormats()
}
